# Refitting NCDEs to individual patients (for manuscript figure):
# adds a ground-truth "Label" column (0 = Control, 1 = MDD) and refreshes
# the per-patient fit statistics from the rerun.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Label" header column (H), matching the bold/centered/bordered style of the other headers
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Per-row Label values (0 = Control, 1 = MDD) and refreshed fit statistics
$ws.Range("D2").Value = 0.5822443915529257
$ws.Range("E2").Value = 0.5822443915529257
$ws.Range("H2").Value = 0

$ws.Range("D3").Value = 0.4201958726310183
$ws.Range("E3").Value = 0.4201958726310183
$ws.Range("H3").Value = 0

$ws.Range("H4").Value = 0

$ws.Range("D5").Value = 0.3449823408148606
$ws.Range("E5").Value = 0.3449823408148606
$ws.Range("H5").Value = 0

$ws.Range("D6").Value = 0.3918419951777649
$ws.Range("E6").Value = 0.3918419951777649
$ws.Range("H6").Value = 0

$ws.Range("H7").Value = 1

$ws.Range("H8").Value = 1

$ws.Range("H9").Value = 1

$ws.Range("D10").Value = 0.4150803942412987
$ws.Range("E10").Value = 0.5849196057587014
$ws.Range("H10").Value = 1

$ws.Range("F11").Value = 292.3453063964844
$ws.Range("H11").Value = 1

$ws.Range("H12").Value = 0

$ws.Range("H13").Value = 0

$ws.Range("H14").Value = 0

$ws.Range("H15").Value = 0

$ws.Range("H16").Value = 0

$ws.Range("H17").Value = 1

$ws.Range("H18").Value = 1

$ws.Range("H19").Value = 1

$ws.Range("H20").Value = 1

$ws.Range("H21").Value = 1

